# final draft text & metadata edits before review
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Correct / update personnel entries (row 7 and row 8 swap + text fixes):
#   Row 7 becomes "S. Alejandra" / "Casillo Cieza"
#   Row 8 becomes "Arshia" / "Mehta"
$ws.Cells.Item(7, 1).Value = "S. Alejandra"
$ws.Cells.Item(7, 3).Value = "Casillo Cieza"
$ws.Cells.Item(8, 1).Value = "Arshia"
$ws.Cells.Item(8, 3).Value = "Mehta"

# Update the active selection to reflect where editing finished
$ws.Range("C8").Select()
